$wb = $excel.ActiveWorkbook

# --- Append Week 15 simulation data to shared-string play-log cells ---
$wb.Worksheets.Item("YDS").Range("B2").Value = "6 1 2 7 0 -2 0 5 1 5 8 9 5 2 2 15 3 3 8 2 2 4 13 2 -4 3 -3 5 1 9 10 6 4 2 4 0 5 2 10 5 11 6 2 -5 1 2 2 6 4 0 4 4 2 5 1 4 -2 0 2 11 2 3 1 4 1 5 15 0 6 6 3 7 15 2 9 15 9 2 0 3 3 6 4 11 3 8 5 7 0 6 0 1 4 3 -3 10 2 6 15 3 6 0 1 3 -1 7 5 29 8 -1 0 2 12 0 2 -1 0 4 1 0 22 9 15 0 4 3 4 2 2 8 0 8 6 10 5 3 1 3 0 1 9 6 5 2 4 5 0 3 3 0 4 3 5 10 0 3 1 1 1 4 0 -2 3 8 4 4 0 15 5 0 16 2 7 19 2 3 4 6 3 -2 4 6 1 10 4 0 5 6 2 14 3 16 5 3 4 1 14 4 12 2 5 2 11 4 6 0 4 3 0 8 0 10 6 1 3 11 4 3 6 7 5 6 6 3 1 2 3 7 5 7 2 0 3 5 3 5 6 0 7 6 0 5 9 2 2 9 4 9 -2 3 0 3 6 8 6 -1 6 7 5 0 25 2 3 2 4 3 19 -1 11 12 4 0 -1 2 5 2 2 0 0 10 3 6 6 0 2"
$wb.Worksheets.Item("YDS").Range("B3").Value = "7 67 1 18 4 17 17 1 6 19 56 6 37 3 6 17 15 9 2 16 23 14 13 16 9 4 8 16 11 14 43 2 11 44 10 14 7 3 7 6 12 22 11 5 6 6 6 9 2 11 21 22 2 75 0 20 5 40 10 5 17 6 4 10 3 20 5 18 5 2 7 14 14 7 35 18 9 0 6 8 7 22 8 6 7 13 9 12 5 9 14 -1 16 11 9 4 17 17 25 5 9 5 8 6 7 9 68 20 20 13 3 23 1 24 33 13 30 16 6 28 10 15 9 4 4 -6 9 -4 16 25 4 4 5 18 5 6 13 16 29 4 5 1 5 5 2 11 11 27 3 22 10 4 9 2 7 13 22 5 14 59 14 5 6 18 5 10 5 8 15 14 3 4 15 9 3 11 28 2 11 52 5 68 11 20 6 5 8 14 12 4 0 7 1 10 6 18 20 26 15 6 17 6 3 4 9 22 13 7 2 13 3 16 5 9 6 9 3 5 2 6 8 16 29 13 10 8 33 5 -2 5 2 6 -3 12 16 6 11 5 9 5 2 25 9 4 5 79 5 11 22 3 19 6 1 8 9 54 13 7 10 7 7 7 11 14 -1 27 4 35 5 9 8 7 7 7 4 19 8 7 43 29 8 3 19 22 2 5 5 3 8 1 5 13 1 7 11 8 2 17 6 40 6 2 6 52 4 4 23 4 17 44 7 12 -4"
$wb.Worksheets.Item("YDS").Range("C2").Value = "41 0 14 2 -3 -2 1 3 5 10 3 3 0 9 3 1 6 4 4 6 1 3 8 5 4 3 3 4 13 0 0 0 5 2 8 4 0 5 12 1 14 6 2 9 -1 1 1 9 3 3 1 4 1 3 1 5 -1 4 3 0 2 3 1 4 9 19 -1 0 10 -2 10 -1 18 8 1 3 4 12 9 3 3 5 0 4 1 1 1 5 4 4 2 1 0 3 54 7 4 14 0 11 3 -3 1 3 13 1 4 4 0 7 4 2 3 4 4 3 6 -4 5 4 2 2 6 2 1 5 12 4 14 9 2 3 0 2 2 0 6 -1 2 4 2 2 5 4 6 1 1 8 3 7 6 3 1 5 12 3 14 1 -3 3 4 0 1 3 2 2 11 14 4 5 0 4 5 2 0 0 2 1 4 2 14 4 1 4 1 3 1 3 -3 8 -7 0 1 4 2 14 2 11 6 3 4 2 1 6 0 3 2 1 4 7 5 1 9 4 0 2 0 8 6 5 -1 9 -1 1 6 2 13 2 8 17 4 -2 4 4 4 -6 3 -3 4 8 4 2 0 -2 2 1 9 2 0 -4 8 3 6 2 4 1 -1 0 1 2 7 5 0 8 4 1 8 2 5 2 4 -4 3 2 5 -1 2 7 7 3 3 2 1 3 1 -1 10 0 2 9 -1 7 1 1 -3 13 -1 1 2 -1 4 0 2 1 5 2 13 1 1 0 0 1 7 2 5 12 8 16 0 17 6 2 1 8 1 1"
$wb.Worksheets.Item("YDS").Range("C3").Value = "4 9 5 10 4 19 -2 10 9 7 3 10 9 8 6 7 10 8 3 10 1 11 5 11 9 9 10 5 6 14 10 9 5 9 42 2 18 7 13 7 3 7 12 23 8 0 15 9 34 10 8 5 1 15 17 4 3 7 7 0 17 5 26 3 20 20 32 14 8 20 15 11 6 7 17 11 -4 23 13 31 17 8 -2 11 8 5 5 21 -1 9 9 7 5 8 25 5 41 6 14 10 7 9 -3 14 7 10 3 16 9 7 9 18 6 12 25 5 25 4 25 14 19 29 9 8 13 5 1 17 8 3 21 23 7 5 14 20 13 17 8 3 16 2 9 4 14 7 8 3 10 10 11 6 10 3 10 1 11 12 11 11 7 6 2 4 63 -1 21 12 15 8 20 14 9 7 5 37 6 6 5 4 3 6 5 17 2 10 6 8 0 9 26 -1 15 15 3 -5 1 6 7 8 8 10 7 6 45 24 8 8 9 9 23 2 8 31 15 5 -3 7 2 10 14 1 8 3 5 9 9 4 7 11 7 13 16 -1 18 6 7 8 10 20 21 6 3 19 7 11 7 8 15 0 2 7 7 5 1 54 27 43 7 10 6 7 8 28 6 6 5 2 9 2 8 6 7 9 3 28 -1 6 11 5 16 9 -2 12 9 7 17 10 9 14 6 17 11 14 5 7 2 18 -2 26 10 41 13 3 9 0 5 12 6 6 7 23 15 36 17 12 8 -2 3 11 9"
$wb.Worksheets.Item("ST").Range("D3").Value = "50 33 46 30 41 51 27 37 52 45 33 54 48 38 34 52 36 41 46 49 48 35 57 52 45 45 38 42 44 53 41 41 33 46 46 58 38 42 40"
$wb.Worksheets.Item("ST").Range("D4").Value = "3 0 0 0 0 4 0 0 13 0 0 0 0 0 0 0 0 0 0 0 1 0 0 11 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0"
$wb.Worksheets.Item("ST").Range("D5").Value = "0 0 0 9 0 0 0 0 4 0 0 11 0 8 28 0 11 0 0 0 2 0 8 0 0 0 0 0 0 0 7 0 0 0 0 0 13 0 0 0 0 0 0 0 0"
$wb.Worksheets.Item("ST").Range("B6").Value = "21 21 28 16 25 20 12 18 19 7 17 17 17 24 8 15 19 26 17 17 16 23 20 28 65 16 27"

# --- Update aggregate numeric cells to reflect the new totals ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 146
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 34
$ws.Range("G2").Value = 41
$ws.Range("J2").Value = 24
$ws.Range("N2").Value = 11
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 162
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 19
$ws.Range("F3").Value = 104
$ws.Range("G3").Value = 41
$ws.Range("I3").Value = 45
$ws.Range("J3").Value = 53
$ws.Range("L3").Value = 243
$ws.Range("M3").Value = 159
$ws.Range("Q3").Value = 422

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 161
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = 6
$ws.Range("G2").Value = 55
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 23
$ws.Range("N2").Value = 20
$ws.Range("O2").Value = 27
$ws.Range("P2").Value = 13
$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 145
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 32
$ws.Range("F3").Value = 101
$ws.Range("G3").Value = 29
$ws.Range("H3").Value = 21
$ws.Range("I3").Value = 58
$ws.Range("J3").Value = 40
$ws.Range("L3").Value = 274
$ws.Range("M3").Value = 176
$ws.Range("Q3").Value = 480

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 78
$ws.Range("D2").Value = 39
$ws.Range("F2").Value = 103
$ws.Range("G2").Value = 97
$ws.Range("J2").Value = 38
$ws.Range("K2").Value = 35
$ws.Range("N2").Value = 12
$ws.Range("O2").Value = 9
$ws.Range("B3").Value = 49

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 8

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("D2").Value = 7
$ws.Range("B3").Value = 8
$ws.Range("B4").Value = 2
